$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 85
$prev = $row - 1

# Replicate the formatting (styles) of the previous data row onto the new row
$ws.Range("A$prev`:V$prev").Copy()
$ws.Range("A$row`:V$row").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 84
$ws.Cells.Item($row, 2).Value = "denmark"
$ws.Cells.Item($row, 3).Value = "superliga"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45236.79166666666
$ws.Cells.Item($row, 6).Value = "Brondby"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Aarhus"
$ws.Cells.Item($row, 9).Value = 1

$ws.Cells.Item($row, 10).Value = 1.98
$ws.Cells.Item($row, 11).Value = "30/10/2023 19:13"
$ws.Cells.Item($row, 12).Value = 2.1
$ws.Cells.Item($row, 13).Value = "06/11/2023 18:50"

$ws.Cells.Item($row, 14).Value = 3.52
$ws.Cells.Item($row, 15).Value = "30/10/2023 19:13"
$ws.Cells.Item($row, 16).Value = 3.32
$ws.Cells.Item($row, 17).Value = "06/11/2023 18:43"

$ws.Cells.Item($row, 18).Value = 3.9
$ws.Cells.Item($row, 19).Value = "30/10/2023 19:13"
$ws.Cells.Item($row, 20).Value = 3.89
$ws.Cells.Item($row, 21).Value = "06/11/2023 18:56"

$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/denmark/superliga/brondby-aarhus/IT01An6G/"
